$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Utilisateurs" to "Courses"
$ws.Name = "Courses"

# Clear existing data (columns A:F, rows 1:2) and remove extra columns D:F
$ws.Cells.Clear()

# Header row
$ws.Range("A1").Value = "id"
$ws.Range("B1").Value = "item"
$ws.Range("C1").Value = "quantity"

# Data rows
$data = @(
    @(1745483727576, "Gruyère (kg)", 1),
    @(1745483739301, "Bolo (pot)", 2),
    @(1745483754454, "Sauce tomate (pot)", 1),
    @(1745483763510, "Pesto verde (pot)", 1),
    @(1745483768824, "Pesto rouge (pot)", 1),
    @(1745483782830, "Beurre (plaquette)", 1),
    @(1745483808182, "Jambon (tranche)", 20),
    @(1745483815286, "Fromage (tranche)", 20),
    @(1745483829826, "Chips (paquet)", 3),
    @(1745483841437, "Nutella (pot)", 1),
    @(1745483857388, "Bière (L)", 10),
    @(1745483866138, "Vodka (L)", 4),
    @(1745483888070, "Smirnoff (Bouteille)", 14),
    @(1745483896109, "Bulles (Bouteille)", 2),
    @(1745483904748, "Ricard (Bouteille)", 1),
    @(1745483912320, "Grenadine (Bouteille)", 1),
    @(1745483920571, "Jaeger (Bouteille)", 2),
    @(1745483927641, "Limonade (Bouteille)", 2),
    @(1745483935346, "Taurine (Bouteille)", 6),
    @(1745486439425, "Multifruit (Bouteille)", 1)
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $row++
}
